# Update "GDP per Capita" values for the Dominican Republic (Data sheet).
# The existing 61 rows (years 1950-2010) get revised values, and 6 new rows
# are appended for years 2011-2016.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# All 67 GDP-per-capita values, in row order, for years 1950..2016.
$gdpValues = @(
    "1637",
    "1780",
    "1868",
    "1790",
    "1833",
    "1886",
    "2008",
    "2066",
    "2104",
    "2047",
    "2075",
    "1964",
    "2222",
    "2291",
    "2365",
    "2007",
    "2206",
    "2211",
    "2150",
    "2316",
    "2488",
    "2678",
    "2928",
    "3196",
    "3295",
    "3365",
    "3504",
    "3586",
    "3583",
    "3649",
    "3781",
    "3846",
    "3826",
    "3905",
    "3830",
    "3653",
    "3682",
    "3877",
    "3821",
    "4219",
    "3939",
    "3969.56867380379",
    "4185.40593822095",
    "4490.01622077702",
    "4597.15081760226",
    "4851.62876552566",
    "5198.56664916418",
    "5616.88992165113",
    "6012.69255026626",
    "6419.78109287858",
    "6788.0838992463",
    "6917.03038105384",
    "7327.36595752144",
    "7324.69734804766",
    "7441.44762923987",
    "8156.4309532226",
    "9059.88442531908",
    "9866.45803030943",
    "10225.5288424606",
    "10365.6748748019",
    "11276.4608119942",
    "11679",
    "11848",
    "12251",
    "13017",
    "13762",
    "14489"
)

$firstDataRow = 2
$firstYear = 1950
$countryCode = 214.0
$countryName = "Dominican Republic"
$indicatorName = "GDP per Capita"

$totalRows = $gdpValues.Length
$lastRow = $firstDataRow + $totalRows - 1

# Make sure columns A-D are populated (existing rows already have the right
# values; this also fills in the 6 brand-new rows at the bottom).
for ($i = 0; $i -lt $totalRows; $i++) {
    $row = $firstDataRow + $i
    $year = $firstYear + $i
    $ws.Cells.Item($row, 1).Value = $countryCode
    $ws.Cells.Item($row, 2).Value = $countryName
    $ws.Cells.Item($row, 3).Value = $indicatorName
    $ws.Cells.Item($row, 4).Value = [double]$year
}

# Column E ("Data") must stay text (as in the source workbook), so route the
# writes through a text formula + copy/paste-values instead of .Value, which
# would silently coerce numeric-looking strings to numbers.
$helperCol = 10
for ($i = 0; $i -lt $totalRows; $i++) {
    $helperRow = $i + 1
    $ws.Cells.Item($helperRow, $helperCol).Formula = "=""" + $gdpValues[$i] + """"
}

$helperRange = $ws.Range($ws.Cells.Item(1, $helperCol), $ws.Cells.Item($totalRows, $helperCol))
$helperRange.Copy()

$targetRange = $ws.Range($ws.Cells.Item($firstDataRow, 5), $ws.Cells.Item($lastRow, 5))
$targetRange.PasteSpecial(-4163)

$helperRange.Clear()
